# "Report formulas & format"
# The single meaningful content change in this revision is that the
# placeholder "." value that used to sit in A3 (a stray default used while
# the report template's header row was being laid out) is removed, leaving
# the cell blank while keeping its existing style/formatting (s="4").
#
# Clearing the cell content (rather than deleting the whole row/cell) also
# causes the now-unused shared string "." to drop out of the shared string
# table on save, which in turn shifts every other shared-string index down
# by one - exactly what the target workbook shows for A1 and the L2:P2
# header cells (their underlying text is unchanged, only the string-table
# index moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "." placeholder from A3 while preserving its formatting.
$ws.Range("A3").ClearContents()
